$wb = $excel.ActiveWorkbook

# --- Sheet "Hoja1": update the daily conversion message in A1 ---
$wsHoja1 = $wb.Worksheets.Item("Hoja1")
$oldText = $wsHoja1.Range("A1").Value()
$newText = $oldText -replace [regex]::Escape("✅ 1000 Bs = 9.79 = 41143.08 pesos`n✅ 41143.08 pesos = 9.78 = 971.25 Bs"), "✅ 1000 Bs = 9.6 = 40249.52 pesos`n✅ 40249.52 pesos = 9.56 = 958.39 Bs"
$wsHoja1.Range("A1").Value = $newText

# --- Sheet "tasas": update the exchange-rate figures ---
$wsTasas = $wb.Worksheets.Item("tasas")
$wsTasas.Range("N10").Value = 104.2
$wsTasas.Range("O10").Value = 4194
$wsTasas.Range("O12").Value = 100.221
